$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($refAddr, $val) {
    # Force the cell to remain plain text even when the string looks
    # like a number (e.g. "243.30"), matching the inlineStr/shared-string
    # cells produced by the source scraper. Reset the style back to
    # "Normal" afterwards so no stray number-format style lingers on the cell.
    $r = $ws.Range($refAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '35.284.61'
$ws.Range('E2').Value = '  +0.46%  '

$ws.Range('D3').Value = '1.895.66'
$ws.Range('E3').Value = '  +2.35%  '

$ws.Range('E4').Value = '  +0.12%  '

Set-TextValue 'D5' '243.30'
$ws.Range('E5').Value = '  +2.22%  '

Set-TextValue 'D6' '0.652'
$ws.Range('E6').Value = '  +5.00%  '

$ws.Range('E7').Value = '  +0.07%  '

Set-TextValue 'D8' '41.48'
$ws.Range('E8').Value = '  -1.87%  '

$ws.Range('E9').Value = '  +4.59%  '

Set-TextValue 'D10' '50.09'
$ws.Range('E10').Value = '  +7.72%  '

Set-TextValue 'D11' '0.0708'
$ws.Range('E11').Value = '  +2.30%  '

Set-TextValue 'D12' '0.0998'
$ws.Range('E12').Value = '  +0.71%  '

$ws.Range('D13').Value = '2.171.83'
$ws.Range('E13').Value = '  +2.45%  '

Set-TextValue 'D14' '12.04'
$ws.Range('E14').Value = '  +5.66%  '

Set-TextValue 'D15' '0.692'

$ws.Range('D16').Value = '1.889.50'
$ws.Range('E16').Value = '  +2.12%  '

$ws.Range('E17').Value = '  +1.20%  '

$ws.Range('D18').Value = '35.299.50'
$ws.Range('E18').Value = '  +0.61%  '

Set-TextValue 'D19' '71.28'
$ws.Range('E19').Value = '  +1.78%  '

$ws.Range('D20').Value = '0.0₃0814'
$ws.Range('E20').Value = '  +2.74%  '

Set-TextValue 'D21' '241.11'
$ws.Range('E21').Value = '  +0.24%  '

Set-TextValue 'D22' '12.49'
$ws.Range('E22').Value = '  +2.75%  '

Set-TextValue 'D23' '4.71'
$ws.Range('E23').Value = '  -1.01%  '

$ws.Range('E24').Value = '  +0.07%  '

Set-TextValue 'D25' '2.39'
$ws.Range('E25').Value = '  +31.17%  '

$ws.Range('E26').Value = '  +1.24%  '

Set-TextValue 'D27' '170.07'
$ws.Range('E27').Value = '  +0.18%  '

Set-TextValue 'D28' '8.32'
$ws.Range('E28').Value = '  +3.85%  '

Set-TextValue 'D29' '18.18'
$ws.Range('E29').Value = '  +3.35%  '

$ws.Range('E30').Value = '  +1.87%  '

$ws.Range('E31').Value = '  +2.77%  '

Set-TextValue 'D32' '0.0560'
$ws.Range('E32').Value = '  +1.30%  '

$ws.Range('E33').Value = '  -0.01%  '

Set-TextValue 'D34' '0.923'
$ws.Range('E34').Value = '  +15.74%  '

$ws.Range('E35').Value = '  +1.76%  '

$ws.Range('E36').Value = '  +0.93%  '

$ws.Range('E37').Value = '  +1.53%  '

Set-TextValue 'D38' '1.33'
$ws.Range('E38').Value = '  +2.15%  '

Set-TextValue 'D39' '0.0208'
$ws.Range('E39').Value = '  +3.47%  '

Set-TextValue 'D40' '1.09'
$ws.Range('E40').Value = '  +1.28%  '

Set-TextValue 'D41' '0.0635'
$ws.Range('E41').Value = '  +14.87%  '

Set-TextValue 'D42' '15.75'
$ws.Range('E42').Value = '  +5.14%  '

Set-TextValue 'D43' '88.99'
$ws.Range('E43').Value = '  -1.46%  '

$ws.Range('D44').Value = '1.337.86'
$ws.Range('E44').Value = '  -0.51%  '

Set-TextValue 'D45' '2.36'
$ws.Range('E45').Value = '  +2.11%  '

Set-TextValue 'D46' '46.89'
$ws.Range('E46').Value = '  +36.43%  '

$ws.Range('E47').Value = '  -1.68%  '

$ws.Range('E48').Value = '  +1.34%  '

$ws.Range('E49').Value = '  +0.16%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.080.09'
$ws.Range('E50').Value = '  +2.31%  '

$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextValue 'D51' '11.31'
$ws.Range('E51').Value = '  -15.22%  '
